$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-01-27 07:46:20"
$wsZh.Range("G4").Value = "2016-01-27 07:47:05"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-01-27 07:46:34"
$wsDe.Range("G4").Value = "2016-01-27 07:47:25"
